$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header in B1 from "Fullname" to "Point"
$ws.Range("B1").Value = "Point"

# Move the active selection to B3 (matches the saved cursor position)
$ws.Range("B3").Select()
